$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are textual (e.g. "27.408.14", "1.000") and must remain
# stored as text, not be auto-converted to numbers by Excel. We temporarily
# force a Text number format, assign the value, then restore the original
# style so the cell formatting is left unchanged.

$origStyle = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.408.14"
$ws.Range("D2").Style = $origStyle
$ws.Range("E2").Value = "  +1.34%  "
$origStyle = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.822.32"
$ws.Range("D3").Style = $origStyle
$ws.Range("E3").Value = "  -0.17%  "
$origStyle = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = $origStyle
$ws.Range("E4").Value = "  -0.02%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.18"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +0.82%  "
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = $origStyle
$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4489"
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = "  +2.02%  "
$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3751"
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = "  +1.92%  "
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07497"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  +3.09%  "
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8876"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  +4.98%  "
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.08"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  +1.87%  "
$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.826.70"
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = "  +0.10%  "
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.757"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = "  +1.40%  "
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "93.93"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  +4.26%  "
$ws.Range("E15").Value = "  +2.01%  "
$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07110"
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("E17").Value = "  -0.06%  "
$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008786"
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("E19").Value = "  +0.01%  "
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.18"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  +1.63%  "
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.403.01"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  +1.26%  "
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.322"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  +3.29%  "
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.94"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  +0.48%  "
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.056.16"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  +0.50%  "
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.959"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  -1.94%  "
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.369"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  +7.50%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("E28").Value = "  +1.54%  "
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.374"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  +2.72%  "
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.11"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  +0.93%  "
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08870"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  +0.81%  "
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7886"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  +6.42%  "
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.202"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  +1.83%  "
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.538"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  +2.52%  "
$ws.Range("E35").Value = "  +1.13%  "
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9999"
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  +1.40%  "
$ws.Range("E38").Value = "  +2.59%  "
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05330"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  +1.80%  "
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.406"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  +2.01%  "
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5326"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  +3.11%  "
$ws.Range("E42").Value = "  +1.46%  "
$ws.Range("E43").Value = "  -0.52%  "
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.305"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  +19.28%  "
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.723"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  +2.18%  "
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5102"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  +5.73%  "
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.61"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  -0.24%  "
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "105.78"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("E49").Value = "  +2.46%  "
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.000"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  +0.04%  "
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06383"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  +0.74%  "
